$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.738254
$ws.Range("H2").Value = 2.214762
$ws.Range("I2").Value = 0.005691320045803731
$ws.Range("J2").Value = 0.005691320045803731
$ws.Range("Q2").Value = 0.2027447273426667
$ws.Range("R2").Value = 1.824702546084
$ws.Range("S2").Value = 0.005691320045803731
$ws.Range("T2").Value = 0.005691320045803731

# Row 3 updates
$ws.Range("I3").Value = 0.9440493064670392
$ws.Range("J3").Value = 0.9440493064670391
$ws.Range("S3").Value = 0.9440493064670392
$ws.Range("T3").Value = 0.9440493064670391

# Row 4 updates
$ws.Range("G4").Value = 5.698467
$ws.Range("H4").Value = 17.095401
$ws.Range("I4").Value = 0.0439304080539368
$ws.Range("J4").Value = 0.04393040805393679
$ws.Range("Q4").Value = 1.564954796298
$ws.Range("R4").Value = 14.084593166682
$ws.Range("S4").Value = 0.0439304080539368
$ws.Range("T4").Value = 0.04393040805393679

# Row 5 updates
$ws.Range("G5").Value = 0.8209666666666666
$ws.Range("H5").Value = 2.4629
$ws.Range("I5").Value = 0.006328965433220369
$ws.Range("J5").Value = 0.006328965433220369
$ws.Range("Q5").Value = 0.2254598864222222
$ws.Range("R5").Value = 2.0291389778
$ws.Range("S5").Value = 0.006328965433220369
$ws.Range("T5").Value = 0.006328965433220369
